$d = $word.ActiveDocument

$d.Content.Find.Execute("59+23=", $true, $false, $false, $false, $false, $true, 1, $false, "96-76=", 2) | Out-Null
$d.Content.Find.Execute("54+45=", $true, $false, $false, $false, $false, $true, 1, $false, "27+13=", 2) | Out-Null
$d.Content.Find.Execute("9+74=", $true, $false, $false, $false, $false, $true, 1, $false, "64-45=", 2) | Out-Null
$d.Content.Find.Execute("75-63=", $true, $false, $false, $false, $false, $true, 1, $false, "8-1=", 2) | Out-Null
$d.Content.Find.Execute("5+47=", $true, $false, $false, $false, $false, $true, 1, $false, "3+0=", 2) | Out-Null
$d.Content.Find.Execute("45-37=", $true, $false, $false, $false, $false, $true, 1, $false, "49-16=", 2) | Out-Null
$d.Content.Find.Execute("67-52=", $true, $false, $false, $false, $false, $true, 1, $false, "60-39=", 2) | Out-Null
$d.Content.Find.Execute("18-15=", $true, $false, $false, $false, $false, $true, 1, $false, "45+41=", 2) | Out-Null
$d.Content.Find.Execute("32-18=", $true, $false, $false, $false, $false, $true, 1, $false, "84+7=", 2) | Out-Null
$d.Content.Find.Execute("17+72=", $true, $false, $false, $false, $false, $true, 1, $false, "18+7=", 2) | Out-Null
$d.Content.Find.Execute("33-11=", $true, $false, $false, $false, $false, $true, 1, $false, "35-3=", 2) | Out-Null
$d.Content.Find.Execute("94-71=", $true, $false, $false, $false, $false, $true, 1, $false, "67-17=", 2) | Out-Null
$d.Content.Find.Execute("51-36=", $true, $false, $false, $false, $false, $true, 1, $false, "62-48=", 2) | Out-Null
$d.Content.Find.Execute("66+8=", $true, $false, $false, $false, $false, $true, 1, $false, "19+47=", 2) | Out-Null
$d.Content.Find.Execute("19+28=", $true, $false, $false, $false, $false, $true, 1, $false, "81-20=", 2) | Out-Null
$d.Content.Find.Execute("6+16=", $true, $false, $false, $false, $false, $true, 1, $false, "26-9=", 2) | Out-Null
$d.Content.Find.Execute("75-61=", $true, $false, $false, $false, $false, $true, 1, $false, "90-1=", 2) | Out-Null
$d.Content.Find.Execute("72-7=", $true, $false, $false, $false, $false, $true, 1, $false, "69+9=", 2) | Out-Null
$d.Content.Find.Execute("7+1=", $true, $false, $false, $false, $false, $true, 1, $false, "64+17=", 2) | Out-Null
$d.Content.Find.Execute("51+24=", $true, $false, $false, $false, $false, $true, 1, $false, "91-64=", 2) | Out-Null
$d.Content.Find.Execute("85-13=", $true, $false, $false, $false, $false, $true, 1, $false, "24-12=", 2) | Out-Null
$d.Content.Find.Execute("36+38=", $true, $false, $false, $false, $false, $true, 1, $false, "4+90=", 2) | Out-Null
$d.Content.Find.Execute("98-71=", $true, $false, $false, $false, $false, $true, 1, $false, "58-24=", 2) | Out-Null
$d.Content.Find.Execute("34-29=", $true, $false, $false, $false, $false, $true, 1, $false, "66-38=", 2) | Out-Null
$d.Content.Find.Execute("25+14=", $true, $false, $false, $false, $false, $true, 1, $false, "14+38=", 2) | Out-Null
$d.Content.Find.Execute("34-28=", $true, $false, $false, $false, $false, $true, 1, $false, "13+34=", 2) | Out-Null
$d.Content.Find.Execute("33+6=", $true, $false, $false, $false, $false, $true, 1, $false, "46+53=", 2) | Out-Null
$d.Content.Find.Execute("27+4=", $true, $false, $false, $false, $false, $true, 1, $false, "16+54=", 2) | Out-Null
$d.Content.Find.Execute("39+49=", $true, $false, $false, $false, $false, $true, 1, $false, "52+45=", 2) | Out-Null
$d.Content.Find.Execute("19+23=", $true, $false, $false, $false, $false, $true, 1, $false, "94-42=", 2) | Out-Null
$d.Content.Find.Execute("9+33=", $true, $false, $false, $false, $false, $true, 1, $false, "54-23=", 2) | Out-Null
$d.Content.Find.Execute("31+52=", $true, $false, $false, $false, $false, $true, 1, $false, "87-58=", 2) | Out-Null
$d.Content.Find.Execute("38+6=", $true, $false, $false, $false, $false, $true, 1, $false, "64+2=", 2) | Out-Null
$d.Content.Find.Execute("72-28=", $true, $false, $false, $false, $false, $true, 1, $false, "67-34=", 2) | Out-Null
$d.Content.Find.Execute("63-36=", $true, $false, $false, $false, $false, $true, 1, $false, "21+17=", 2) | Out-Null
$d.Content.Find.Execute("76-11=", $true, $false, $false, $false, $false, $true, 1, $false, "43-7=", 2) | Out-Null
$d.Content.Find.Execute("0+48=", $true, $false, $false, $false, $false, $true, 1, $false, "61-54=", 2) | Out-Null
$d.Content.Find.Execute("58+21=", $true, $false, $false, $false, $false, $true, 1, $false, "94-76=", 2) | Out-Null
$d.Content.Find.Execute("49-5=", $true, $false, $false, $false, $false, $true, 1, $false, "58+7=", 2) | Out-Null
$d.Content.Find.Execute("81+16=", $true, $false, $false, $false, $false, $true, 1, $false, "52-25=", 2) | Out-Null
$d.Content.Find.Execute("82-51=", $true, $false, $false, $false, $false, $true, 1, $false, "41-36=", 2) | Out-Null
$d.Content.Find.Execute("61-52=", $true, $false, $false, $false, $false, $true, 1, $false, "79-44=", 2) | Out-Null
$d.Content.Find.Execute("56-46=", $true, $false, $false, $false, $false, $true, 1, $false, "11-5=", 2) | Out-Null
$d.Content.Find.Execute("19+29=", $true, $false, $false, $false, $false, $true, 1, $false, "75-16=", 2) | Out-Null
$d.Content.Find.Execute("61-43=", $true, $false, $false, $false, $false, $true, 1, $false, "42+14=", 2) | Out-Null
$d.Content.Find.Execute("78-46=", $true, $false, $false, $false, $false, $true, 1, $false, "23+74=", 2) | Out-Null
$d.Content.Find.Execute("69-16=", $true, $false, $false, $false, $false, $true, 1, $false, "28+9=", 2) | Out-Null
$d.Content.Find.Execute("47+17=", $true, $false, $false, $false, $false, $true, 1, $false, "64-37=", 2) | Out-Null
$d.Content.Find.Execute("64+9=", $true, $false, $false, $false, $false, $true, 1, $false, "58-40=", 2) | Out-Null
$d.Content.Find.Execute("1+57=", $true, $false, $false, $false, $false, $true, 1, $false, "37-33=", 2) | Out-Null
$d.Content.Find.Execute("73-64=", $true, $false, $false, $false, $false, $true, 1, $false, "28-16=", 2) | Out-Null
$d.Content.Find.Execute("38+51=", $true, $false, $false, $false, $false, $true, 1, $false, "66-6=", 2) | Out-Null
$d.Content.Find.Execute("0+79=", $true, $false, $false, $false, $false, $true, 1, $false, "46+3=", 2) | Out-Null
$d.Content.Find.Execute("84-39=", $true, $false, $false, $false, $false, $true, 1, $false, "35-2=", 2) | Out-Null
$d.Content.Find.Execute("0+63=", $true, $false, $false, $false, $false, $true, 1, $false, "80+9=", 2) | Out-Null
$d.Content.Find.Execute("4+95=", $true, $false, $false, $false, $false, $true, 1, $false, "45+44=", 2) | Out-Null
$d.Content.Find.Execute("1+27=", $true, $false, $false, $false, $false, $true, 1, $false, "56+42=", 2) | Out-Null
$d.Content.Find.Execute("77+0=", $true, $false, $false, $false, $false, $true, 1, $false, "84-67=", 2) | Out-Null
$d.Content.Find.Execute("21+14=", $true, $false, $false, $false, $false, $true, 1, $false, "72+19=", 2) | Out-Null
$d.Content.Find.Execute("31+22=", $true, $false, $false, $false, $false, $true, 1, $false, "57+0=", 2) | Out-Null
$d.Content.Find.Execute("87-81=", $true, $false, $false, $false, $false, $true, 1, $false, "40+45=", 2) | Out-Null
$d.Content.Find.Execute("20+73=", $true, $false, $false, $false, $false, $true, 1, $false, "96-89=", 2) | Out-Null
$d.Content.Find.Execute("52-35=", $true, $false, $false, $false, $false, $true, 1, $false, "65-39=", 2) | Out-Null
$d.Content.Find.Execute("71-70=", $true, $false, $false, $false, $false, $true, 1, $false, "24-3=", 2) | Out-Null
$d.Content.Find.Execute("0+65=", $true, $false, $false, $false, $false, $true, 1, $false, "40-19=", 2) | Out-Null
$d.Content.Find.Execute("57+19=", $true, $false, $false, $false, $false, $true, 1, $false, "59+18=", 2) | Out-Null
$d.Content.Find.Execute("6+29=", $true, $false, $false, $false, $false, $true, 1, $false, "91-31=", 2) | Out-Null
$d.Content.Find.Execute("61-2=", $true, $false, $false, $false, $false, $true, 1, $false, "97-31=", 2) | Out-Null
$d.Content.Find.Execute("42+2=", $true, $false, $false, $false, $false, $true, 1, $false, "30-26=", 2) | Out-Null
$d.Content.Find.Execute("21+16=", $true, $false, $false, $false, $false, $true, 1, $false, "49+36=", 2) | Out-Null
$d.Content.Find.Execute("71-47=", $true, $false, $false, $false, $false, $true, 1, $false, "45-16=", 2) | Out-Null
$d.Content.Find.Execute("46+27=", $true, $false, $false, $false, $false, $true, 1, $false, "59-1=", 2) | Out-Null
$d.Content.Find.Execute("77+14=", $true, $false, $false, $false, $false, $true, 1, $false, "59-11=", 2) | Out-Null
$d.Content.Find.Execute("95-87=", $true, $false, $false, $false, $false, $true, 1, $false, "47-34=", 2) | Out-Null
$d.Content.Find.Execute("13+24=", $true, $false, $false, $false, $false, $true, 1, $false, "22+67=", 2) | Out-Null
$d.Content.Find.Execute("70-62=", $true, $false, $false, $false, $false, $true, 1, $false, "75-35=", 2) | Out-Null
$d.Content.Find.Execute("35+39=", $true, $false, $false, $false, $false, $true, 1, $false, "49-11=", 2) | Out-Null
$d.Content.Find.Execute("85-26=", $true, $false, $false, $false, $false, $true, 1, $false, "17+60=", 2) | Out-Null
$d.Content.Find.Execute("92-29=", $true, $false, $false, $false, $false, $true, 1, $false, "84+13=", 2) | Out-Null
$d.Content.Find.Execute("68-58=", $true, $false, $false, $false, $false, $true, 1, $false, "14+74=", 2) | Out-Null
$d.Content.Find.Execute("89-6=", $true, $false, $false, $false, $false, $true, 1, $false, "44-31=", 2) | Out-Null
$d.Content.Find.Execute("5+51=", $true, $false, $false, $false, $false, $true, 1, $false, "6+15=", 2) | Out-Null
$d.Content.Find.Execute("64-42=", $true, $false, $false, $false, $false, $true, 1, $false, "15+63=", 2) | Out-Null
$d.Content.Find.Execute("11-4=", $true, $false, $false, $false, $false, $true, 1, $false, "95-61=", 2) | Out-Null
$d.Content.Find.Execute("32-21=", $true, $false, $false, $false, $false, $true, 1, $false, "40+14=", 2) | Out-Null
$d.Content.Find.Execute("3+9=", $true, $false, $false, $false, $false, $true, 1, $false, "19+28=", 2) | Out-Null
$d.Content.Find.Execute("66-14=", $true, $false, $false, $false, $false, $true, 1, $false, "39+10=", 2) | Out-Null
$d.Content.Find.Execute("35-20=", $true, $false, $false, $false, $false, $true, 1, $false, "47-13=", 2) | Out-Null
$d.Content.Find.Execute("76+5=", $true, $false, $false, $false, $false, $true, 1, $false, "41+19=", 2) | Out-Null
$d.Content.Find.Execute("54+36=", $true, $false, $false, $false, $false, $true, 1, $false, "8+39=", 2) | Out-Null
$d.Content.Find.Execute("40-11=", $true, $false, $false, $false, $false, $true, 1, $false, "43-17=", 2) | Out-Null
$d.Content.Find.Execute("38+10=", $true, $false, $false, $false, $false, $true, 1, $false, "30+23=", 2) | Out-Null
$d.Content.Find.Execute("6+52=", $true, $false, $false, $false, $false, $true, 1, $false, "30-11=", 2) | Out-Null
$d.Content.Find.Execute("39+4=", $true, $false, $false, $false, $false, $true, 1, $false, "63-59=", 2) | Out-Null
$d.Content.Find.Execute("35-27=", $true, $false, $false, $false, $false, $true, 1, $false, "45+49=", 2) | Out-Null
$d.Content.Find.Execute("59-18=", $true, $false, $false, $false, $false, $true, 1, $false, "3+60=", 2) | Out-Null
$d.Content.Find.Execute("82-68=", $true, $false, $false, $false, $false, $true, 1, $false, "49-36=", 2) | Out-Null
$d.Content.Find.Execute("3+65=", $true, $false, $false, $false, $false, $true, 1, $false, "46+11=", 2) | Out-Null
$d.Content.Find.Execute("24-24=", $true, $false, $false, $false, $false, $true, 1, $false, "82+5=", 2) | Out-Null
$d.Content.Find.Execute("79+15=", $true, $false, $false, $false, $false, $true, 1, $false, "73-2=", 2) | Out-Null
